$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (keeps them as text, matching source data)
$textCells = @("D5", "D6", "D10", "D11", "D15", "D17", "D21", "D22", "D23", "D24", "D27", "D29", "D30", "D31", "D34", "D36", "D41", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "42.143.66"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3
$ws.Range("D3").Value = "2.263.52"
$ws.Range("E3").Value = "  -1.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "307.63"
$ws.Range("E5").Value = "  +0.34%  "

# Row 6
$ws.Range("D6").Value = "97.13"
$ws.Range("E6").Value = "  +0.61%  "

# Row 7
$ws.Range("E7").Value = "  -1.63%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  -1.55%  "

# Row 10
$ws.Range("D10").Value = "35.07"
$ws.Range("E10").Value = "  -3.76%  "

# Row 11
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -2.15%  "

# Row 12
$ws.Range("E12").Value = "  +0.48%  "

# Row 13
$ws.Range("E13").Value = "  +1.06%  "

# Row 14
$ws.Range("D14").Value = "2.611.72"
$ws.Range("E14").Value = "  -1.12%  "

# Row 15
$ws.Range("D15").Value = "14.58"
$ws.Range("E15").Value = "  -0.14%  "

# Row 16
$ws.Range("D16").Value = "2.244.79"
$ws.Range("E16").Value = "  -0.92%  "

# Row 17
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").Value = "  -2.01%  "

# Row 18
$ws.Range("D18").Value = "41.943.80"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19
$ws.Range("E19").Value = "  -4.97%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  -2.05%  "

# Row 21
$ws.Range("D21").Value = "5.96"

# Row 22
$ws.Range("D22").Value = "67.54"
$ws.Range("E22").Value = "  -0.57%  "

# Row 23
$ws.Range("D23").Value = "236.27"
$ws.Range("E23").Value = "  -2.83%  "

# Row 24
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +1.05%  "

# Row 25
$ws.Range("E25").Value = "  -0.96%  "

# Row 26
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("D27").Value = "23.50"
$ws.Range("E27").Value = "  -2.18%  "

# Row 28
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "9.51"
$ws.Range("E29").Value = "  -1.02%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +0.56%  "

# Row 31
$ws.Range("D31").Value = "164.44"
$ws.Range("E31").Value = "  +1.65%  "

# Row 32
$ws.Range("E32").Value = "  -2.06%  "

# Row 33
$ws.Range("E33").Value = "  +0.04%  "

# Row 34
$ws.Range("D34").Value = "3.14"
$ws.Range("E34").Value = "  +1.06%  "

# Row 35
$ws.Range("E35").Value = "  -2.36%  "

# Row 36
$ws.Range("D36").Value = "17.44"
$ws.Range("E36").Value = "  +0.15%  "

# Row 37
$ws.Range("E37").Value = "  +0.33%  "

# Row 38
$ws.Range("E38").Value = "  -4.57%  "

# Row 39
$ws.Range("E39").Value = "  -1.23%  "

# Row 40
$ws.Range("E40").Value = "  -3.54%  "

# Row 41
$ws.Range("D41").Value = "4.13"
$ws.Range("E41").Value = "  -1.48%  "

# Row 42
$ws.Range("E42").Value = "  -2.82%  "

# Row 43
$ws.Range("D43").Value = "1.951.70"
$ws.Range("E43").Value = "  -2.71%  "

# Row 44
$ws.Range("D44").Value = "19.03"
$ws.Range("E44").Value = "  -0.82%  "

# Row 45
$ws.Range("D45").Value = "0.0280"
$ws.Range("E45").Value = "  -1.76%  "

# Row 46
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  -3.18%  "

# Row 47
$ws.Range("D47").Value = "9.78"
$ws.Range("E47").Value = "  -4.15%  "

# Row 48
$ws.Range("D48").Value = "53.43"
$ws.Range("E48").Value = "  -1.16%  "

# Row 49
$ws.Range("D49").Value = "2.483.98"
$ws.Range("E49").Value = "  -1.11%  "

# Row 50
$ws.Range("D50").Value = "92.36"
$ws.Range("E50").Value = "  +0.30%  "

# Row 51
$ws.Range("D51").Value = "71.56"
$ws.Range("E51").Value = "  -1.68%  "

# Restore default style on cells where we forced text number format
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
